$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 84.2
$ws.Range("D2").Value = 97.40000000000001
$ws.Range("E2").Value = 90.8
$ws.Range("F2").Value = 93.5
$ws.Range("G2").Value = 77.3
$ws.Range("H2").Value = 89.3
$ws.Range("I2").Value = 86.90000000000001
$ws.Range("J2").Value = 83.90000000000001
$ws.Range("C4").Value = 450.7
$ws.Range("D4").Value = 494.1
$ws.Range("E4").Value = 437
$ws.Range("F4").Value = 457.2
$ws.Range("G4").Value = 403
$ws.Range("H4").Value = 491.2
$ws.Range("I4").Value = 373.1
$ws.Range("J4").Value = 407
$ws.Range("C6").Value = 64.7
$ws.Range("D6").Value = 76.09999999999999
$ws.Range("E6").Value = 68.40000000000001
$ws.Range("F6").Value = 70
$ws.Range("H6").Value = 77.5
$ws.Range("I6").Value = 65.7
$ws.Range("C8").Value = 29.6
$ws.Range("D8").Value = 35.6
$ws.Range("E8").Value = 31
$ws.Range("F8").Value = 36.6
$ws.Range("G8").Value = 35.8
$ws.Range("H8").Value = 45.9
$ws.Range("I8").Value = 35.5
$ws.Range("C10").Value = 70.2
$ws.Range("D10").Value = 82.09999999999999
$ws.Range("F10").Value = 75.59999999999999
$ws.Range("G10").Value = 64.90000000000001
$ws.Range("H10").Value = 79.90000000000001
$ws.Range("I10").Value = 71.59999999999999
$ws.Range("C12").Value = 470.5
$ws.Range("D12").Value = 525.2
$ws.Range("E12").Value = 480.5
$ws.Range("F12").Value = 489.3
$ws.Range("G12").Value = 443
$ws.Range("H12").Value = 542.4
$ws.Range("I12").Value = 446.2
$ws.Range("J12").Value = 459.6
$ws.Range("C14").Value = 157.4
$ws.Range("D14").Value = 169.5
$ws.Range("E14").Value = 153.8
$ws.Range("F14").Value = 170.2
$ws.Range("G14").Value = 161.2
$ws.Range("H14").Value = 162.8
$ws.Range("I14").Value = 160.4
$ws.Range("J14").Value = 160.8
$ws.Range("C16").Value = 101.2
$ws.Range("D16").Value = 116
$ws.Range("E16").Value = 103.5
$ws.Range("F16").Value = 107.7
$ws.Range("H16").Value = 120.5
$ws.Range("C18").Value = 7165.4
$ws.Range("D18").Value = 7832.2
$ws.Range("E18").Value = 7209.6
$ws.Range("F18").Value = 7681.7
$ws.Range("G18").Value = 6930.8
$ws.Range("H18").Value = 7730.5
$ws.Range("I18").Value = 7072.2
$ws.Range("J18").Value = 7326
$ws.Range("C20").Value = 4119.8
$ws.Range("D20").Value = 4425.5
$ws.Range("E20").Value = 3961.6
$ws.Range("F20").Value = 4285.6
$ws.Range("G20").Value = 3911.2
$ws.Range("H20").Value = 4273
$ws.Range("I20").Value = 3836.6
$ws.Range("J20").Value = 3996.6
$ws.Range("C22").Value = 124.3
$ws.Range("D22").Value = 147.6
$ws.Range("E22").Value = 131.9
$ws.Range("F22").Value = 136.4
$ws.Range("G22").Value = 111.5
$ws.Range("H22").Value = 131.4
$ws.Range("I22").Value = 114.4
$ws.Range("J22").Value = 130.2
$ws.Range("C24").Value = 236.2
$ws.Range("D24").Value = 267.4
$ws.Range("E24").Value = 233.7
$ws.Range("F24").Value = 261.9
$ws.Range("G24").Value = 193.1
$ws.Range("H24").Value = 247.4
$ws.Range("I24").Value = 209.1
$ws.Range("J24").Value = 216.5
$ws.Range("C26").Value = 83.09999999999999
$ws.Range("D26").Value = 88.8
$ws.Range("E26").Value = 78.09999999999999
$ws.Range("F26").Value = 91.3
$ws.Range("G26").Value = 73.59999999999999
$ws.Range("H26").Value = 85.7
$ws.Range("I26").Value = 77.3
$ws.Range("C28").Value = 195.1
$ws.Range("E28").Value = 192.2
$ws.Range("F28").Value = 200.1
$ws.Range("G28").Value = 162
$ws.Range("H28").Value = 182.2
$ws.Range("I28").Value = 149.8
$ws.Range("J28").Value = 174.7
$ws.Range("C30").Value = 341.6
$ws.Range("E30").Value = 341.1
$ws.Range("F30").Value = 397.7
$ws.Range("G30").Value = 317
$ws.Range("I30").Value = 324
$ws.Range("J30").Value = 342.1
$ws.Range("C32").Value = 115.2
$ws.Range("D32").Value = 112.1
$ws.Range("E32").Value = 110.4
$ws.Range("F32").Value = 118.8
$ws.Range("G32").Value = 111.5
$ws.Range("H32").Value = 123.6
$ws.Range("I32").Value = 98.8
$ws.Range("C34").Value = 91.2
$ws.Range("D34").Value = 101.5
$ws.Range("E34").Value = 97.40000000000001
$ws.Range("F34").Value = 108.7
$ws.Range("G34").Value = 93.09999999999999
$ws.Range("H34").Value = 110.2
$ws.Range("J34").Value = 100.9
$ws.Range("C36").Value = 140.5
$ws.Range("D36").Value = 153.9
$ws.Range("E36").Value = 129.5
$ws.Range("F36").Value = 144.5
$ws.Range("G36").Value = 117.6
$ws.Range("H36").Value = 133.4
$ws.Range("I36").Value = 126.6
$ws.Range("C38").Value = 30.3
$ws.Range("D38").Value = 38.3
$ws.Range("E38").Value = 29.6
$ws.Range("F38").Value = 34.8
$ws.Range("G38").Value = 37.4
$ws.Range("I38").Value = 36.6
$ws.Range("C40").Value = 170.3
$ws.Range("D40").Value = 208.7
$ws.Range("E40").Value = 187.4
$ws.Range("F40").Value = 199.2
$ws.Range("G40").Value = 182
$ws.Range("H40").Value = 210.5
$ws.Range("I40").Value = 193.4
$ws.Range("J40").Value = 195.1
$ws.Range("C42").Value = 48.5
$ws.Range("D42").Value = 49.7
$ws.Range("E42").Value = 49.8
$ws.Range("H42").Value = 58
$ws.Range("I42").Value = 50.9
$ws.Range("C44").Value = 25.1
$ws.Range("D44").Value = 25.5
$ws.Range("E44").Value = 23.3
$ws.Range("F44").Value = 28.9
$ws.Range("H44").Value = 22.6
$ws.Range("I44").Value = 18.9
$ws.Range("J44").Value = 21
$ws.Range("C46").Value = 399.9
$ws.Range("D46").Value = 447.7
$ws.Range("E46").Value = 415.5
$ws.Range("F46").Value = 415.8
$ws.Range("G46").Value = 395.6
$ws.Range("H46").Value = 443.5
$ws.Range("I46").Value = 377.5
$ws.Range("D48").Value = 270
$ws.Range("E48").Value = 239
$ws.Range("F48").Value = 236.5
$ws.Range("G48").Value = 224.4
$ws.Range("H48").Value = 241.7
$ws.Range("I48").Value = 234.6
$ws.Range("C50").Value = 184.7
$ws.Range("D50").Value = 202.9
$ws.Range("E50").Value = 191.6
$ws.Range("F50").Value = 199.8
$ws.Range("G50").Value = 190
$ws.Range("H50").Value = 214.5
$ws.Range("I50").Value = 158
$ws.Range("C52").Value = 79.90000000000001
$ws.Range("D52").Value = 93.09999999999999
$ws.Range("F52").Value = 97.7
$ws.Range("G52").Value = 79.7
$ws.Range("H52").Value = 86.5
$ws.Range("I52").Value = 78.8
$ws.Range("J52").Value = 86.2
$ws.Range("C54").Value = 63.6
$ws.Range("D54").Value = 74
$ws.Range("E54").Value = 67.5
$ws.Range("F54").Value = 72.7
$ws.Range("G54").Value = 66.59999999999999
$ws.Range("H54").Value = 81.09999999999999
$ws.Range("I54").Value = 67.2
$ws.Range("C56").Value = 185.6
$ws.Range("D56").Value = 204
$ws.Range("E56").Value = 190.5
$ws.Range("F56").Value = 210.9
$ws.Range("G56").Value = 174.6
$ws.Range("H56").Value = 200.8
$ws.Range("I56").Value = 176.4
$ws.Range("D58").Value = 41.5
$ws.Range("E58").Value = 32.1
$ws.Range("F58").Value = 37.9
$ws.Range("G58").Value = 31.1
$ws.Range("H58").Value = 37.1
$ws.Range("I58").Value = 32.1
$ws.Range("J58").Value = 34.8
$ws.Range("C60").Value = 157.5
$ws.Range("D60").Value = 191
$ws.Range("E60").Value = 172.4
$ws.Range("G60").Value = 151.9
$ws.Range("H60").Value = 166.6
$ws.Range("I60").Value = 159.4
$ws.Range("C62").Value = 323.1
$ws.Range("D62").Value = 374.6
$ws.Range("E62").Value = 355.5
$ws.Range("F62").Value = 366.6
$ws.Range("G62").Value = 335.1
$ws.Range("H62").Value = 375.6
$ws.Range("I62").Value = 333.6
$ws.Range("C64").Value = 14.8
$ws.Range("D64").Value = 13.2
$ws.Range("E64").Value = 13.4
$ws.Range("F64").Value = 11.8
$ws.Range("H64").Value = 18.7
$ws.Range("C66").Value = 17.3
$ws.Range("E66").Value = 15.7
$ws.Range("G66").Value = 25.7
$ws.Range("I66").Value = 20.1
